$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds date serial values. Rows 2-9 are all being
# bumped from 45207 (2023-10-08) to 45208 (2023-10-09).
$ws.Range("C2:C9").Value = 45208
